$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = 16.4098
$ws.Range("D7").Value = -7.265200000000005
$ws.Range("C9").Value = -10.23850000000001
$ws.Range("D12").Value = -6.837899999999999
$ws.Range("D14").Value = -7.575500000000003
$ws.Range("E15").Value = 15.9584
$ws.Range("C18").Value = -12.5921
$ws.Range("C20").Value = -11.15630000000001
$ws.Range("D26").Value = -8.4725
$ws.Range("C27").Value = -12.41229999999999
$ws.Range("D27").Value = -8.726699999999999
$ws.Range("D29").Value = -7.345099999999998
$ws.Range("E33").Value = 17.06340000000001
$ws.Range("C35").Value = -11.37870000000001
$ws.Range("E35").Value = 16.65190000000002
$ws.Range("D37").Value = -7.629699999999996
$ws.Range("D38").Value = -7.720599999999997
$ws.Range("E38").Value = 16.62460000000002
$ws.Range("E43").Value = 17.2182
$ws.Range("E44").Value = 16.14139999999999
$ws.Range("E47").Value = 16.4258
$ws.Range("D51").Value = -8.469600000000003
$ws.Range("E51").Value = 16.4182
$ws.Range("D52").Value = -7.816099999999999
$ws.Range("D55").Value = -8.504899999999994
$ws.Range("E57").Value = 16.54960000000001
$ws.Range("E63").Value = 18.18060000000002
$ws.Range("C69").Value = -10.5435
$ws.Range("D69").Value = -7.186399999999995
$ws.Range("D70").Value = -7.491100000000001
$ws.Range("E70").Value = 17.08160000000001
$ws.Range("C76").Value = -12.72910000000001
$ws.Range("C78").Value = -11.13860000000001
$ws.Range("D81").Value = -7.717700000000007
$ws.Range("C82").Value = -12.4275
$ws.Range("C83").Value = -13.96410000000001
$ws.Range("D83").Value = -8.599600000000001
$ws.Range("E88").Value = 16.45600000000001
$ws.Range("C93").Value = -11.1869
$ws.Range("E99").Value = 16.7723
$ws.Range("D102").Value = -7.624399999999997
